$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.419.30"
$ws.Cells.Item(2, 5).Value = "  +0.29%  "
$ws.Cells.Item(3, 4).Value = "'1.701.07"
$ws.Cells.Item(3, 5).Value = "  +0.84%  "
$ws.Cells.Item(4, 4).Value = "'1.009"
$ws.Cells.Item(4, 5).Value = "  +0.09%  "
$ws.Cells.Item(5, 5).Value = "  +0.36%  "
$ws.Cells.Item(6, 4).Value = "'0.5508"
$ws.Cells.Item(6, 5).Value = "  +4.79%  "
$ws.Cells.Item(7, 4).Value = "'1.009"
$ws.Cells.Item(7, 5).Value = "  +0.01%  "
$ws.Cells.Item(8, 4).Value = "'0.2752"
$ws.Cells.Item(9, 4).Value = "'0.06476"
$ws.Cells.Item(9, 5).Value = "  +0.91%  "
$ws.Cells.Item(10, 4).Value = "'22.09"
$ws.Cells.Item(10, 5).Value = "  +0.18%  "
$ws.Cells.Item(11, 4).Value = "'0.07678"
$ws.Cells.Item(11, 5).Value = "  +2.41%  "
$ws.Cells.Item(12, 4).Value = "'1.703.48"
$ws.Cells.Item(12, 5).Value = "  +0.15%  "
$ws.Cells.Item(13, 4).Value = "'4.555"
$ws.Cells.Item(13, 5).Value = "  -0.39%  "
$ws.Cells.Item(14, 4).Value = "'0.5849"
$ws.Cells.Item(14, 5).Value = "  +0.45%  "
$ws.Cells.Item(15, 4).Value = "'0.000008390"
$ws.Cells.Item(15, 5).Value = "  -1.32%  "
$ws.Cells.Item(16, 4).Value = "'65.66"
$ws.Cells.Item(16, 5).Value = "  +1.78%  "
$ws.Cells.Item(17, 4).Value = "'26.455.80"
$ws.Cells.Item(17, 5).Value = "  +0.26%  "
$ws.Cells.Item(18, 4).Value = "'4.952"
$ws.Cells.Item(18, 5).Value = "  +0.33%  "
$ws.Cells.Item(19, 4).Value = "'1.009"
$ws.Cells.Item(19, 5).Value = "  +0.18%  "
$ws.Cells.Item(20, 4).Value = "'10.99"
$ws.Cells.Item(20, 5).Value = "  +0.96%  "
$ws.Cells.Item(21, 4).Value = "'192.54"
$ws.Cells.Item(21, 5).Value = "  +1.56%  "
$ws.Cells.Item(22, 4).Value = "'6.267"
$ws.Cells.Item(22, 5).Value = "  +0.88%  "
$ws.Cells.Item(23, 4).Value = "'1.010"
$ws.Cells.Item(23, 5).Value = "  +0.15%  "
$ws.Cells.Item(24, 4).Value = "'148.88"
$ws.Cells.Item(24, 5).Value = "  +2.95%  "
$ws.Cells.Item(25, 4).Value = "'0.1330"
$ws.Cells.Item(25, 5).Value = "  +7.91%  "
$ws.Cells.Item(26, 4).Value = "'7.934"
$ws.Cells.Item(26, 5).Value = "  +2.71%  "
$ws.Cells.Item(27, 4).Value = "'15.84"
$ws.Cells.Item(27, 5).Value = "  -0.02%  "
$ws.Cells.Item(28, 4).Value = "'0.06312"
$ws.Cells.Item(28, 5).Value = "  -5.16%  "
$ws.Cells.Item(29, 4).Value = "'1.380"
$ws.Cells.Item(29, 5).Value = "  +2.15%  "
$ws.Cells.Item(30, 5).Value = "  +0.14%  "
$ws.Cells.Item(31, 5).Value = "  +0.92%  "
$ws.Cells.Item(32, 4).Value = "'3.615"
$ws.Cells.Item(32, 5).Value = "  +1.18%  "
$ws.Cells.Item(33, 4).Value = "'1.688"
$ws.Cells.Item(33, 5).Value = "  +1.61%  "
$ws.Cells.Item(34, 4).Value = "'1.046"
$ws.Cells.Item(34, 5).Value = "  +1.78%  "
$ws.Cells.Item(35, 4).Value = "'0.6177"
$ws.Cells.Item(35, 5).Value = "  -0.94%  "
$ws.Cells.Item(36, 4).Value = "'2.411"
$ws.Cells.Item(36, 5).Value = "  +0.47%  "
$ws.Cells.Item(37, 4).Value = "'2.722"
$ws.Cells.Item(37, 5).Value = "  +0.71%  "
$ws.Cells.Item(38, 4).Value = "'0.01660"
$ws.Cells.Item(38, 5).Value = "  +2.48%  "
$ws.Cells.Item(39, 4).Value = "'6.188"
$ws.Cells.Item(39, 5).Value = "  -3.04%  "
$ws.Cells.Item(40, 4).Value = "'1.119.54"
$ws.Cells.Item(40, 5).Value = "  +0.58%  "
$ws.Cells.Item(41, 4).Value = "'0.8834"
$ws.Cells.Item(41, 5).Value = "  -0.09%  "
$ws.Cells.Item(42, 5).Value = "  -0.24%  "
$ws.Cells.Item(43, 4).Value = "'101.44"
$ws.Cells.Item(43, 5).Value = "  +0.34%  "
$ws.Cells.Item(44, 4).Value = "'1.848.75"
$ws.Cells.Item(44, 5).Value = "  +0.63%  "
$ws.Cells.Item(45, 4).Value = "'57.82"
$ws.Cells.Item(45, 5).Value = "  +1.48%  "
$ws.Cells.Item(46, 5).Value = "  -3.05%  "
$ws.Cells.Item(47, 4).Value = "'8.235"
$ws.Cells.Item(47, 5).Value = "  +0.63%  "
$ws.Cells.Item(48, 4).Value = "'1.008"
$ws.Cells.Item(48, 5).Value = "  -0.15%  "
$ws.Cells.Item(49, 5).Value = "  +0.15%  "
$ws.Cells.Item(50, 4).Value = "'6.125"
$ws.Cells.Item(50, 5).Value = "  +0.80%  "
$ws.Cells.Item(51, 4).Value = "'0.4304"
$ws.Cells.Item(51, 5).Value = "  -0.06%  "
